$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("PlotSettings")
$ws.Rows.Item(1).Delete()
